# Commit: "removing restriction on dhw types"
#
# The HOT_WATER assembly sheet gains a new "class_dhw" column (inserted
# between the existing "code" and "Tsww0_C" columns) that tags each DHW
# assembly with its temperature class (HIGH_TEMP / MEDIUM_TEMP / LOW_TEMP),
# replacing what used to be an implicit restriction baked into the row
# order/description alone.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HOT_WATER")

# Insert a new column at C - this shifts the old C (Tsww0_C) and D
# (Qwwmax_Wm2) columns right to D and E, carrying their data/format/
# comments with them, and copies the left-neighbour (B) formatting into
# the freshly inserted column C.
$ws.Columns("C").Insert()

# Populate the new column's data cells first (top-to-bottom by value
# groups) so new shared strings are interned in the same order the
# original author entered them, then set the header last.
$ws.Range("C3").Value = "HIGH_TEMP"
$ws.Range("C5").Value = "LOW_TEMP"
$ws.Range("C4").Value = "MEDIUM_TEMP"
$ws.Range("C1").Value = "class_dhw"
$ws.Range("C2").Value = "NONE"
$ws.Range("C6").Value = "HIGH_TEMP"

# Page setup - the workbook now carries explicit print settings for this
# sheet.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Leave the selection the way the author left it: editing column C on the
# HOT_WATER tab, with C2:C6 highlighted and that tab active; COOLING's
# lingering selection moves off of its old C2 spot too.
$wsCooling = $wb.Worksheets.Item("COOLING")
$wsCooling.Range("C1").Select()

$ws.Activate()
$ws.Range("C2:C6").Select()
